$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of kaspa buy data for 2025-06-27. Column A holds the date as the
# literal text "06/27/2025" (matching the other recent rows in the sheet,
# which store dates as plain text rather than date serials). Pre-format the
# cell as Text before writing so Excel doesn't auto-convert the MM/DD/YYYY
# string into a date value, then clear the formatting back off so the cell
# is left with the default (unstyled) look, same as its neighbours.
$dateCell = $ws.Range("A36")
$dateCell.NumberFormat = "@"
$dateCell.Value = "06/27/2025"
$dateCell.ClearFormats()

$ws.Range("B36").Value = 674.4099999999962
$ws.Range("C36").Value = 0.07413887694429247
$ws.Range("D36").Value = 50
